$wb = $excel.ActiveWorkbook

$oldId = "2423d5c1-52c7-428d-bbbb-0b8a74148bf3"
$newId = "7366215b-7069-49b9-a0d9-0cf019ca4725"

$oldMdName = "$oldId.md"
$newMdName = "$newId.md"

$oldZhHash = "80364fa12a0d31fc6c17c26c8cfcb8c224e45caa"
$newZhHash = "2d37b9f4c6b149cb1a265f7f5ee16c06fab6352e"
$oldDeHash = "80364fa12a0d31fc6c17c26c8cfcb8c224e45caa"
$newDeHash = "2d37b9f4c6b149cb1a265f7f5ee16c06fab6352e"

$oldZhXlfName = "$oldId.$oldZhHash.zh-cn.xlf"
$newZhXlfName = "$newId.$newZhHash.zh-cn.xlf"
$oldDeXlfName = "$oldId.$oldDeHash.de-de.xlf"
$newDeXlfName = "$newId.$newDeHash.de-de.xlf"

$newZhDatetime = "2016-03-03 07:58:24"
$newDeDatetime = "2016-03-03 07:58:35"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/8795060ad17c3a9a864aa8e993f421e6b7a419f5/e2e/$oldMdName"
$configAddress = "https://github.com/OpenLocalizationTest/oltest/blob/8795060ad17c3a9a864aa8e993f421e6b7a419f5/.localization-config"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/127536147edf4d6a7871334f4e4e13a09bf51a20/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$oldZhXlfName"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/856cc54c457a71fed9fcfcd1e149b30f2b3c39c3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$oldDeXlfName"

# ---- Sheet 1: Overview ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdAddress, "", "", $newMdName)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $configAddress, "", "", ".localization-config")

# ---- Sheet 2: zh-cn ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("D2").Value = $newZhDatetime
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdAddress, "", "", $newMdName)
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zhXlfAddress, "", "", $newZhXlfName)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $configAddress, "", "", ".localization-config")

# ---- Sheet 3: de-de ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("D2").Value = $newDeDatetime
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdAddress, "", "", $newMdName)
$ws3.Hyperlinks.Add($ws3.Range("C2"), $deXlfAddress, "", "", $newDeXlfName)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $configAddress, "", "", ".localization-config")

$wb.Save()
